$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

Write-Output ("RowsBefore=" + $table.Rows.Count)

# Remove 5 of the 6 trailing empty rows, keeping exactly one empty row
# after the last data row.
for ($i = 0; $i -lt 5; $i++) {
    $table.Rows.Item($table.Rows.Count).Delete()
}

Write-Output ("RowsAfter=" + $table.Rows.Count)
